$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear out the old data rows (rows 2 and 3, columns A..S) so leftover
# values/formatting don't linger once we shrink/replace them.
$ws.Range("A2:S3").Clear()

# New single-column list of preview-sound identifiers replacing the old
# multi-column row layout.
$ws.Range("A2").Value = "tontwi"
$ws.Range("A3").Value = "ohdsub"
$ws.Range("A4").Value = "mksskn"
$ws.Range("A5").Value = "sw2op"
$ws.Range("A6").Value = "vfgbs"
$ws.Range("A7").Value = "ohdmi2"

# Update the selection to match the saved view.
$ws.Range("A13").Select()
